# Updated latest Guinea master data.
# Re-shapes the device_type sheet: reorders the original 5 columns
# (code, name, descr, lang_code, is_active) and appends the standard
# MOSIP master-data audit columns (cr_by, cr_dtimes, upd_by, upd_dtimes,
# is_deleted, del_dtimes), populating every data row (2-8) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$headers = @("code", "name", "descr", "lang_code", "is_active", "cr_by", "cr_dtimes", "upd_by", "upd_dtimes", "is_deleted", "del_dtimes")
for ($j = 0; $j -lt $headers.Length; $j++) {
    $ws.Cells.Item(1, $j + 1).Value = $headers[$j]
}

# ---- Data rows (code, name, descr, lang_code) ------------------------
$rows = @(
    @("FRS", "Scanner d'empreintes digitales", "Pour capturer les empreintes digitales", "fra"),
    @("IRS", "Scanner d'iris", "Pour capturer l'iris", "fra"),
    @("CMR", "Appareil photo", "Pour prendre la photo du visage", "fra"),
    @("SCN", "Scanner de documents", "Pour numÃ©riser les documents", "fra"),
    @("PRT", "Imprimante", "Pour imprimer des documents", "fra"),
    @("PTS", "Imprimante Scanner", "Pour imprimer et scanner les documents", "fra"),
    @("QRS", "Scanner QR code", "Pour lire les QR codes", "fra")
)

# Created timestamp shared by every row (2023-06-02 ~13:50 UTC serialized
# as an Excel date serial), formatted like a built-in mm:ss.0 code (47)
# to match how the source workbook displays the cr_dtimes column.
$crDtimes = 45079.576914178244
$ws.Range("G2:G8").NumberFormat = "mm:ss.0"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]   # code
    $ws.Cells.Item($r, 2).Value = $row[1]   # name
    $ws.Cells.Item($r, 3).Value = $row[2]   # descr
    $ws.Cells.Item($r, 4).Value = $row[3]   # lang_code
    $ws.Cells.Item($r, 5).Value = $true     # is_active

    $ws.Cells.Item($r, 6).Value = "superadmin"   # cr_by
    $ws.Cells.Item($r, 7).Value = $crDtimes       # cr_dtimes
    $ws.Cells.Item($r, 8).Value = "NULL"          # upd_by
    $ws.Cells.Item($r, 9).Value = "NULL"          # upd_dtimes
    $ws.Cells.Item($r, 10).Value = $false         # is_deleted
    $ws.Cells.Item($r, 11).Value = "NULL"         # del_dtimes
}

# ---- Sheet view / selection tweaks -----------------------------------
$excel.Goto($ws.Range("D13"))
